$wb = $excel.ActiveWorkbook
$win = $wb.Windows.Item(1)
$win.Width = 26420
$win.Height = 13440
Write-Host "set done"
Write-Host ("Win width: " + $win.Width)
Write-Host ("Win height: " + $win.Height)
